# Applies the edits described by the diff:
#  - various text replacements inside the two header/body tables
#  - removal of the "Url del recurso" and "Usuario" rows from the first table
#  - rename of the remaining "USB" value to "URL"
#  - author name / phone / email updates in the second table

$d = $word.ActiveDocument

# --- Table 1 --------------------------------------------------------------
$t1 = $d.Tables.Item(1)

$t1.Cell(2, 4).Range.Text  = "Prueba para rdd Video"
$t1.Cell(4, 3).Range.Text  = "Medio superior"
$t1.Cell(4, 6).Range.Text  = "Duración en semanas"
$t1.Cell(4, 10).Range.Text = "4"
$t1.Cell(5, 3).Range.Text  = 'Centro de Estudios Científicos y Tecnológicos No. 2 "Miguel Bernard"'
$t1.Cell(5, 11).Range.Text = "2/9/2024"
$t1.Cell(6, 4).Range.Text  = "Dibujo técnico"
$t1.Cell(7, 4).Range.Text  = "Mixta"
$t1.Cell(8, 2).Range.Text  = "Dibujo del rostro humano"
$t1.Cell(9, 4).Range.Text  = "URL"

# Remove the two trailing rows ("Url del recurso" / "Usuario") entirely.
# Delete from the bottom up so row indices stay valid.
$t1.Rows.Item(11).Delete()
$t1.Rows.Item(10).Delete()

# --- Table 2 (autores) -----------------------------------------------------
$t2 = $d.Tables.Item(2)

$t2.Cell(2, 7).Range.Text  = "Ricardo Gutiérrez Florez"
$t2.Cell(3, 4).Range.Text  = "57426"
$t2.Cell(3, 10).Range.Text = "cgonzalp@ipn.mx"
